$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade row (row 8) appended to the repeater output.
# Copy formatting from the previous row's date / boolean cells first so the
# new cells reuse the existing style (date number format on A and G)
# instead of Excel minting a brand-new style entry.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("G7").Copy()
$ws.Range("G8").PasteSpecial(-4122)

$ws.Range("A8").Value = 42654.746458333335
$ws.Range("B8").Value = $true
$ws.Range("C8").Value = 9980.33
$ws.Range("D8").Value = 9965.3799999999992
$ws.Range("E8").Value = 104.43
$ws.Range("F8").Value = 104.74
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = 0.3
$ws.Range("I8").Value = $false

# Column A needs to widen slightly to fit the new (longer) date value.
$ws.Columns.Item(1).ColumnWidth = 14.5
